$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "en"
$ws.Cells.Item(2, 2).Value = 15
$ws.Cells.Item(2, 3).Value = "Iva852"

# Row 3
$ws.Cells.Item(3, 1).Value = "de"
$ws.Cells.Item(3, 2).Value = 7
$ws.Cells.Item(3, 3).Value = "formschub"

# Row 4
$ws.Cells.Item(4, 1).Value = "sv"
$ws.Cells.Item(4, 2).Value = 1
$ws.Cells.Item(4, 3).Value = "PGKurki"

# Row 5
$ws.Cells.Item(5, 1).Value = "ja"
$ws.Cells.Item(5, 2).Value = 10
$ws.Cells.Item(5, 3).Value = "Guren"

# Row 6
$ws.Cells.Item(6, 1).Value = "pt"
$ws.Cells.Item(6, 2).Value = 1
$ws.Cells.Item(6, 3).Value = "rcdc3"

# Row 7
$ws.Cells.Item(7, 1).Value = "nl"
$ws.Cells.Item(7, 2).Value = 1
$ws.Cells.Item(7, 3).Value = "debedachtzamen"

# Row 8
$ws.Cells.Item(8, 1).Value = "zh"
$ws.Cells.Item(8, 2).Value = 1
$ws.Cells.Item(8, 3).Value = "chronos"

# Row 9
$ws.Cells.Item(9, 1).Value = "zh-TW"
$ws.Cells.Item(9, 2).Value = 1
$ws.Cells.Item(9, 3).Value = "cxyozn"

# Row 10
$ws.Cells.Item(10, 1).Value = "fr"
$ws.Cells.Item(10, 2).Value = 1
$ws.Cells.Item(10, 3).Value = "Koopa"

# Row 11
$ws.Cells.Item(11, 1).Value = ""
$ws.Cells.Item(11, 2).Value = 2
$ws.Cells.Item(11, 3).Value = "HMJN"
